# Generate Report for Handback
#
# The localization-status report is re-sorted alphabetically by source file
# name, and the row(s) whose handback just completed (257d61c9...md and its
# dependent 5138aaef...md) move from "Ready for handoff" to
# "Handed back: in sync with en-US" with fresh handback timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$ov.Range("A4").Value = "5138aaef-fea0-4fa3-addc-f35200ff812b.md"
$ov.Range("B4").Value = "Handed back: in sync with en-US"
$ov.Range("C4").Value = "Handed back: in sync with en-US"

$ov.Range("A5").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$ov.Range("B5").Value = "Handback transform failed"
$ov.Range("C5").Value = "Handback transform failed"

$ov.Range("A6").Value = ".localization-config"
$ov.Range("B6").Value = "Not to be localized"
$ov.Range("C6").Value = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-09 03:20:44"
$zh.Range("E2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$zh.Range("F2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-09 03:21:45"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-09 03:15:56"
$zh.Range("E3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$zh.Range("F3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-09 03:16:53"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = "5138aaef-fea0-4fa3-addc-f35200ff812b.md"
$zh.Range("B4").Value = "Handed back: in sync with en-US"
$zh.Range("C4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf"
$zh.Range("D4").Value = "2016-03-09 03:20:44"
$zh.Range("E4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$zh.Range("F4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf"
$zh.Range("G4").Value = "2016-03-09 03:21:45"
$zh.Range("H4").Value = "Include"

$zh.Range("A5").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$zh.Range("B5").Value = "Handback transform failed"
$zh.Range("C5").Value = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf"
$zh.Range("D5").Value = "2016-03-09 03:18:30"
$zh.Range("E5").Value = ""
$zh.Range("F5").Value = ""
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Include"

$zh.Range("A6").Value = ".localization-config"
$zh.Range("B6").Value = "Not to be localized"
$zh.Range("C6").Value = ""
$zh.Range("D6").Value = "0001-01-01 00:00:00"
$zh.Range("E6").Value = ""
$zh.Range("F6").Value = ""
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf"
$de.Range("D2").Value = "2016-03-09 03:20:54"
$de.Range("E2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$de.Range("F2").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf"
$de.Range("G2").Value = "2016-03-09 03:22:20"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf"
$de.Range("D3").Value = "2016-03-09 03:16:11"
$de.Range("E3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$de.Range("F3").Value = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf"
$de.Range("G3").Value = "2016-03-09 03:17:29"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = "5138aaef-fea0-4fa3-addc-f35200ff812b.md"
$de.Range("B4").Value = "Handed back: in sync with en-US"
$de.Range("C4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf"
$de.Range("D4").Value = "2016-03-09 03:20:54"
$de.Range("E4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.md"
$de.Range("F4").Value = "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf"
$de.Range("G4").Value = "2016-03-09 03:22:20"
$de.Range("H4").Value = "Include"

$de.Range("A5").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$de.Range("B5").Value = "Handback transform failed"
$de.Range("C5").Value = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf"
$de.Range("D5").Value = "2016-03-09 03:18:41"
$de.Range("E5").Value = ""
$de.Range("F5").Value = ""
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Include"

$de.Range("A6").Value = ".localization-config"
$de.Range("B6").Value = "Not to be localized"
$de.Range("C6").Value = ""
$de.Range("D6").Value = "0001-01-01 00:00:00"
$de.Range("E6").Value = ""
$de.Range("F6").Value = ""
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Ignored"
